$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.305.81"
$ws.Range("E2").Value = "  -5.91%  "
$ws.Range("D3").Value = "3.298.59"
$ws.Range("E3").Value = "  -5.09%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'561.14"
$ws.Range("E5").Value = "  -4.12%  "
$ws.Range("D6").Value = "'128.90"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.299.06"
$ws.Range("E8").Value = "  -5.07%  "
$ws.Range("D9").Value = "'0.472"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("D10").Value = "'7.36"
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("E11").Value = "  -5.07%  "
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("D13").Value = "3.867.50"
$ws.Range("E13").Value = "  -4.84%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "3.299.71"
$ws.Range("E15").Value = "  -5.02%  "
$ws.Range("E16").Value = "  -5.99%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "60.548.53"
$ws.Range("E17").Value = "  -5.49%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'24.17"
$ws.Range("E18").Value = "  -4.25%  "
$ws.Range("D19").Value = "'5.64"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "'13.32"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "'8.95"
$ws.Range("E21").Value = "  -10.32%  "
$ws.Range("D22").Value = "'350.78"
$ws.Range("E22").Value = "  -8.99%  "
$ws.Range("D23").Value = "'0.553"
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "3.434.55"
$ws.Range("E25").Value = "  -4.95%  "
$ws.Range("D26").Value = "'69.23"
$ws.Range("E26").Value = "  -7.15%  "
$ws.Range("D27").Value = "'0.0000108"
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'7.29"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'7.82"
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("E32").Value = "  -5.78%  "
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "3.333.65"
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("D36").Value = "'22.67"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'6.76"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("D40").Value = "'157.90"
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'40.96"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'4.34"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "'0.740"
$ws.Range("E45").Value = "  -7.11%  "
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'22.73"
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.54"
$ws.Range("E48").Value = "  -4.70%  "
$ws.Range("D49").Value = "'6.68"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "'21.56"
$ws.Range("E50").Value = "  +4.62%  "
$ws.Range("D51").Value = "'0.859"
$ws.Range("E51").Value = "  -4.82%  "
